$d = $word.ActiveDocument
$d.Content.Find.Execute("Todas las rutas (excepto /login) requieren un token JWT.", $true, $false, $false, $false, $false, $true, 1, $false, "Todas las rutas (excepto /login y el /get de platos) requieren un token JWT.", 2)
